$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 52: date column (A) must stay text "05-11-2025" without triggering
# Excel's automatic date conversion / new number-format style, and the
# price column (B) gets the new price string.
$ws.Range("A52").Formula = '="05-11-2025"'
$ws.Range("A52").Copy()
$ws.Range("A52").PasteSpecial(-4163)

$priceText = "The price of gold in India today is " + [char]0x20B9 + "12,148 per gram for 24 karat gold, " + [char]0x20B9 + "11,135 per gram for 22 karat gold and " + [char]0x20B9 + "9,111 per gram for 18 karat gold (also called 999 gold)."
$ws.Range("B52").Value2 = $priceText

$excel.CutCopyMode = 0
